# "Updating the sheet with mesh names"
#
# Adds ten new mesh file-name rows (mesh_in_011 .. mesh_in_020) into column F
# of rows 2-11, matching the already-present "mesh_in_001..010" entries that
# live in F12:F21. The new cells get the same base format as the existing
# F-column cells (left/center aligned, thin border) plus a light highlight
# fill so the freshly-added rows stand out, then the selection is moved to
# H12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Pick up the existing "mesh_in_0xx" cell formatting (font/border/alignment)
# from F12 and apply it to the new F2:F11 range before writing the values.
$ws.Range("F12").Copy()
$ws.Range("F2:F11").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Highlight the newly added cells with a light accent fill so they're
# distinguishable from the pre-existing mesh_in_001..010 list.
$ws.Range("F2:F11").Interior.Color = 15064278   # RGB(214,220,229) ~ Text2, Lighter 80%

# Write the ten new mesh names.
$ws.Range("F2").Value = "mesh_in_011"
$ws.Range("F3").Value = "mesh_in_012"
$ws.Range("F4").Value = "mesh_in_013"
$ws.Range("F5").Value = "mesh_in_014"
$ws.Range("F6").Value = "mesh_in_015"
$ws.Range("F7").Value = "mesh_in_016"
$ws.Range("F8").Value = "mesh_in_017"
$ws.Range("F9").Value = "mesh_in_018"
$ws.Range("F10").Value = "mesh_in_019"
$ws.Range("F11").Value = "mesh_in_020"

# Match the author's final cursor position.
$ws.Range("H12").Select() | Out-Null
